$wb = $excel.ActiveWorkbook

# --- Withdraw History: Amount, Time, Date, Location/Place ---
$ws = $wb.Worksheets.Item("Withdraw History")
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Location/Place"
$ws.Range("E1").Clear()
$ws.Range("F1").Clear()

# --- Deposit History: Amount, Time, Date, Location/Place ---
$ws = $wb.Worksheets.Item("Deposit History")
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Location/Place"
$ws.Range("E1").Clear()
$ws.Range("F1").Clear()

# --- Transfer History: Amount, Time, Date, Person ---
$ws = $wb.Worksheets.Item("Transfer History")
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Person"
$ws.Range("E1").Clear()
$ws.Range("F1").Clear()

# --- Absolute History: Amount, Time, Date, Location/Place/Person ---
$ws = $wb.Worksheets.Item("Absolute History")
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Location/Place/Person"
$ws.Range("E1").Clear()
$ws.Range("F1").Clear()

# Make "Absolute History" the active/selected sheet (was "Amount" before)
$ws.Activate()
